# Weekly update: insert 2 new rows of fresh data (Peru origin, week of 2021-09-09)
# above the existing "1a nueva(o)" / "2a nueva(o)" Camote rows, pushing the
# previous data down by two rows (338-343 -> 340-345).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("338:339").Insert()

# New row 338
$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 44448
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = 100112045
$ws.Range("G338").Value = "Zapallo"
$ws.Range("H338").Value = "Camote"
$ws.Range("I338").Value = "1a nueva(o)"
$ws.Range("J338").Value = 740
$ws.Range("K338").Value = 950
$ws.Range("L338").Value = 1000
$ws.Range("M338").Value = 975
$ws.Range("N338").Value = "$/kilo (volumen en unidades)"
$ws.Range("O338").Value = "Perú"
$ws.Range("P338").Value = 975
$ws.Range("Q338").Value = 1
$ws.Range("R338").Value = "Hortaliza"

# New row 339
$ws.Range("A339").Value = 8
$ws.Range("B339").Value = "Terminal La Palmera de La Serena"
$ws.Range("C339").Value = "Coquimbo"
$ws.Range("D339").Value = 44448
$ws.Range("E339").Value = 4
$ws.Range("F339").Value = 100112045
$ws.Range("G339").Value = "Zapallo"
$ws.Range("H339").Value = "Camote"
$ws.Range("I339").Value = "2a nueva(o)"
$ws.Range("J339").Value = 480
$ws.Range("K339").Value = 850
$ws.Range("L339").Value = 900
$ws.Range("M339").Value = 875
$ws.Range("N339").Value = "$/kilo (volumen en unidades)"
$ws.Range("O339").Value = "Perú"
$ws.Range("P339").Value = 875
$ws.Range("Q339").Value = 1
$ws.Range("R339").Value = "Hortaliza"

# Ensure date cells keep the date/time number format used throughout column D
$ws.Range("D338:D339").NumberFormat = $ws.Range("D340").NumberFormat()
